$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.114.58'
$ws.Range('E2').Value = '  -4.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.689.49'
$ws.Range('E3').Value = '  -7.55%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '494.49'
$ws.Range('E5').Value = '  -6.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.58'
$ws.Range('E6').Value = '  -3.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.522'
$ws.Range('E8').Value = '  -5.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.705.22'
$ws.Range('E9').Value = '  -6.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.85'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('E11').Value = '  -6.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.336'
$ws.Range('E12').Value = '  -4.41%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.156.10'
$ws.Range('E14').Value = '  -7.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.152.36'
$ws.Range('E15').Value = '  -4.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.86'
$ws.Range('E16').Value = '  -7.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.693.91'
$ws.Range('E17').Value = '  -7.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000132'
$ws.Range('E18').Value = '  -6.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.58'
$ws.Range('E19').Value = '  -6.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.64'
$ws.Range('E20').Value = '  -7.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '333.43'
$ws.Range('E21').Value = '  -7.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('E22').Value = '  -7.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '61.53'
$ws.Range('E25').Value = '  -2.91%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.169'
$ws.Range('E26').Value = '  -4.66%  '
$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.414'
$ws.Range('E27').Value = '  -7.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0807'
$ws.Range('E29').Value = '  -6.08%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.24'
$ws.Range('E30').Value = '  -5.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  -5.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.62'
$ws.Range('E33').Value = '  -5.24%  '
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.03'
$ws.Range('E35').Value = '  -6.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.17'
$ws.Range('E36').Value = '  -6.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.895'
$ws.Range('E37').Value = '  -10.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.10'
$ws.Range('E38').Value = '  -8.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.91'
$ws.Range('E39').Value = '  -5.25%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.134.02'
$ws.Range('E41').Value = '  -8.20%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.34'
$ws.Range('E42').Value = '  -8.17%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.41'
$ws.Range('E43').Value = '  -6.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0541'
$ws.Range('E44').Value = '  -4.58%  '
$ws.Range('E45').Value = '  -7.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.33'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '18.30'
$ws.Range('E47').Value = '  -11.65%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.54'
$ws.Range('E48').Value = '  -5.99%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0219'
$ws.Range('E49').Value = '  -5.57%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0871'
$ws.Range('E50').Value = '  -5.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.03'
$ws.Range('E51').Value = '  -6.53%  '
